$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "593.80")
# must be forced to remain TEXT (matching the source inlineStr cells), since
# plain Range.Value assignment auto-converts number-looking strings to numeric
# cells (normal Excel "type into a cell" behavior). We temporarily mark the
# cell as Text-formatted, write the value, then restore the Normal style so the
# cell keeps the same (default) style it had before - just like the target file.

$ws.Range('D2').Value = '63.504.60'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.601.57'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.80'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.60'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +2.01%  '
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.62'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = '3.069.74'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '63.348.49'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('E16').Value = '  +7.36%  '
$ws.Range('D17').Value = '2.626.12'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.46'
$ws.Range('E18').Value = '  +7.74%  '
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.10'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.93'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.56'
$ws.Range('E23').Value = '  +2.22%  '
$ws.Range('E24').Value = '  +4.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.32'
$ws.Range('E25').Value = '  +3.21%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '568.05'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.03'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.19'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.415'
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.62'
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '168.91'
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.62'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.96'
$ws.Range('E43').Value = '  +5.51%  '
$ws.Range('E44').Value = '  +4.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.37'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0254'
$ws.Range('E47').Value = '  +4.38%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.05'
$ws.Range('E48').Value = '  +5.62%  '
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.27'
$ws.Range('E50').Value = '  +3.25%  '
$ws.Range('D51').Value = '0.0₆0231'
$ws.Range('E51').Value = '  +17.72%  '

# Restore default (Normal) style on the text-forced cells so no stray
# number-format style is left attached to them.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
